# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column (B), shifting ASIN..is_holiday_week right by one
#  - rewrite the "Week" labels from zero-padded (W01) to unpadded (W1) form
#  - populate the new Week_Start_Date column with the weekly start dates
#  - re-type the (now shifted) is_holiday_week column as boolean

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new blank column before column B (ASIN etc. shift from B..I to C..J)
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Per-week data: week label, start date, holiday flag
$weeks = @(
    @{ Row = 2;  Label = "W1";  Date = "2025-01-05"; Holiday = $false },
    @{ Row = 3;  Label = "W2";  Date = "2025-01-12"; Holiday = $false },
    @{ Row = 4;  Label = "W3";  Date = "2025-01-19"; Holiday = $false },
    @{ Row = 5;  Label = "W4";  Date = "2025-01-26"; Holiday = $false },
    @{ Row = 6;  Label = "W5";  Date = "2025-02-02"; Holiday = $false },
    @{ Row = 7;  Label = "W6";  Date = "2025-02-09"; Holiday = $false },
    @{ Row = 8;  Label = "W7";  Date = "2025-02-16"; Holiday = $false },
    @{ Row = 9;  Label = "W8";  Date = "2025-02-23"; Holiday = $false },
    @{ Row = 10; Label = "W9";  Date = "2025-03-02"; Holiday = $false },
    @{ Row = 11; Label = "W10"; Date = "2025-03-09"; Holiday = $false },
    @{ Row = 12; Label = "W11"; Date = "2025-03-16"; Holiday = $false },
    @{ Row = 13; Label = "W12"; Date = "2025-03-23"; Holiday = $false },
    @{ Row = 14; Label = "W13"; Date = "2025-03-30"; Holiday = $false },
    @{ Row = 15; Label = "W14"; Date = "2025-04-06"; Holiday = $false },
    @{ Row = 16; Label = "W15"; Date = "2025-04-13"; Holiday = $false },
    @{ Row = 17; Label = "W16"; Date = "2025-04-20"; Holiday = $false }
)

foreach ($wk in $weeks) {
    $r = $wk.Row
    # Column A: unpadded week label (W01 -> W1, etc.)
    $ws.Cells.Item($r, 1).Value = $wk.Label
    # Column B: newly inserted Week_Start_Date column.
    # Prefix with an apostrophe so Excel stores the literal text
    # "2025-01-05" instead of auto-converting it to a date serial number.
    $ws.Cells.Item($r, 2).Value = "'" + $wk.Date
    # Column J (was I, shifted by the column insert): is_holiday_week as a true boolean
    $ws.Cells.Item($r, 10).Value = $wk.Holiday
}
